# Apply "Add data for 2021-12-31" update:
# - rename sheet / update period label from 12-22 to 12-23
# - update December (through 12-22/23) row (row 14) values
# - update Total row (row 15) values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet to reflect the new "through" date
$ws.Name = "Through 2021-12-23"

# Update the row label
$ws.Range("A14").Value = "December (through 12-23)"

# Row 14 - December (through 12-23)
$ws.Range("C14").Value = 28
$ws.Range("D14").Value = 0.125

$ws.Range("F14").Value = 67
$ws.Range("G14").Value = 0.0822

$ws.Range("I14").Value = 82
$ws.Range("J14").Value = 0.1087

$ws.Range("K14").Value = 5
$ws.Range("L14").Value = 48
$ws.Range("M14").Value = 0.0943

$ws.Range("O14").Value = 43
$ws.Range("P14").Value = 0.0851

$ws.Range("Q14").Value = 8
$ws.Range("R14").Value = 104
$ws.Range("S14").Value = 0.0714

$ws.Range("U14").Value = 150
$ws.Range("V14").Value = 0.0132

# Row 15 - Total
$ws.Range("C15").Value = 286
$ws.Range("D15").Value = 0.1146

$ws.Range("F15").Value = 571
$ws.Range("G15").Value = 0.1022

$ws.Range("I15").Value = 840
$ws.Range("J15").Value = 0.08

$ws.Range("K15").Value = 79
$ws.Range("L15").Value = 656
$ws.Range("M15").Value = 0.1075

$ws.Range("O15").Value = 523
$ws.Range("P15").Value = 0.0998

$ws.Range("Q15").Value = 72
$ws.Range("R15").Value = 1304
$ws.Range("S15").Value = 0.0523

$ws.Range("U15").Value = 1693
$ws.Range("V15").Value = 0.0568
